$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3181.721
$ws.Range("I15").Value = 3181.721
$ws.Range("K15").Value = 9545.163
$ws.Range("M15").Value = -9376.163
$ws.Range("H28").Value = 2114.55
$ws.Range("I28").Value = 685.73334
$ws.Range("K28").Value = 685.73334
$ws.Range("M28").Value = -200.73334
$ws.Range("H62").Value = 6955.3125
$ws.Range("I62").Value = 7535.143
$ws.Range("K62").Value = 7535.143
$ws.Range("M62").Value = -6911.143
$ws.Range("H65").Value = 6955.3125
$ws.Range("I65").Value = 7535.143
$ws.Range("K65").Value = 37675.715
$ws.Range("M65").Value = -34555.715
$ws.Range("H86").Value = 166669420
$ws.Range("I86").Value = 250001390
$ws.Range("J86").Value = 5500
$ws.Range("K86").Value = 250001390
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = -250000267
$ws.Range("N86").Value = -7746
$ws.Range("H89").Value = 166669420
$ws.Range("I89").Value = 250001390
$ws.Range("J89").Value = 5500
$ws.Range("K89").Value = 1250006950
$ws.Range("L89").Value = 27500
$ws.Range("M89").Value = -1250001334
$ws.Range("N89").Value = -38732
$ws.Range("H98").Value = 2467.3416
$ws.Range("I98").Value = 2479.025
$ws.Range("K98").Value = 2479.025
$ws.Range("M98").Value = -981.0250000000001
$ws.Range("H106").Value = 2814.8667
$ws.Range("I106").Value = 2730.818
$ws.Range("J106").Value = 3046
$ws.Range("K106").Value = 2730.818
$ws.Range("L106").Value = 3046
$ws.Range("M106").Value = -2099.818
$ws.Range("N106").Value = -4308
$ws.Range("H122").Value = 2467.3416
$ws.Range("I122").Value = 2479.025
$ws.Range("K122").Value = 7437.075000000001
$ws.Range("M122").Value = -4987.075000000001
$ws.Range("H132").Value = 5825
$ws.Range("I132").Value = 5957.0312
$ws.Range("J132").Value = 1600
$ws.Range("K132").Value = 17871.0936
$ws.Range("L132").Value = 4800
$ws.Range("M132").Value = -15341.0936
$ws.Range("N132").Value = -9860
$ws.Range("H135").Value = 695.8461
$ws.Range("I135").Value = 334.25
$ws.Range("J135").Value = 5035
$ws.Range("K135").Value = 3008.25
$ws.Range("L135").Value = 45315
$ws.Range("M135").Value = -473.25
$ws.Range("N135").Value = -50385
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 4760.793
$ws.Range("I137").Value = 1379.4117
$ws.Range("K137").Value = 4138.2351
$ws.Range("M137").Value = -1588.2351

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4783.877
$ws.Range("I32").Value = 3947.4656
$ws.Range("K32").Value = 3947.4656
$ws.Range("M32").Value = -3660.4656
$ws.Range("H56").Value = 23441.666
$ws.Range("J56").Value = 22110
$ws.Range("L56").Value = 22110
$ws.Range("N56").Value = -23594
$ws.Range("H74").Value = 196676.52
$ws.Range("I74").Value = 253741.19
$ws.Range("K74").Value = 253741.19
$ws.Range("M74").Value = -252867.19
$ws.Range("H77").Value = 196676.52
$ws.Range("I77").Value = 253741.19
$ws.Range("K77").Value = 1268705.95
$ws.Range("M77").Value = -1264337.95
$ws.Range("H102").Value = 5252.478
$ws.Range("I102").Value = 4853.4707
$ws.Range("K102").Value = 4853.4707
$ws.Range("M102").Value = -3231.4707
$ws.Range("H132").Value = 4130.75
$ws.Range("I132").Value = 2799.5
$ws.Range("K132").Value = 8398.5
$ws.Range("M132").Value = -5868.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 26322494
$ws.Range("I20").Value = 29418612
$ws.Range("J20").Value = 5499.5
$ws.Range("K20").Value = 29418612
$ws.Range("L20").Value = 5499.5
$ws.Range("M20").Value = -29418365
$ws.Range("N20").Value = -5993.5
$ws.Range("H86").Value = 2910.52
$ws.Range("I86").Value = 2651.5789
$ws.Range("K86").Value = 2651.5789
$ws.Range("M86").Value = -1528.5789
$ws.Range("H89").Value = 2910.52
$ws.Range("I89").Value = 2651.5789
$ws.Range("K89").Value = 13257.8945
$ws.Range("M89").Value = -7641.8945

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1326.0769
$ws.Range("I22").Value = 1014.5
$ws.Range("J22").Value = 2364.6667
$ws.Range("K22").Value = 1014.5
$ws.Range("L22").Value = 2364.6667
$ws.Range("M22").Value = -664.5
$ws.Range("N22").Value = -3064.6667
$ws.Range("H31").Value = 4143.676
$ws.Range("I31").Value = 3152.0386
$ws.Range("J31").Value = 6487.5454
$ws.Range("K31").Value = 3152.0386
$ws.Range("L31").Value = 6487.5454
$ws.Range("M31").Value = -2857.0386
$ws.Range("N31").Value = -7077.5454
$ws.Range("H34").Value = 4143.676
$ws.Range("I34").Value = 3152.0386
$ws.Range("J34").Value = 6487.5454
$ws.Range("K34").Value = 3152.0386
$ws.Range("L34").Value = 6487.5454
$ws.Range("M34").Value = -2950.0386
$ws.Range("N34").Value = -6891.5454
$ws.Range("H105").Value = 2685.5557
$ws.Range("I105").Value = 962.5
$ws.Range("K105").Value = 962.5
$ws.Range("M105").Value = 784.5
$ws.Range("H134").Value = 6565.1177
$ws.Range("I134").Value = 6400.5
$ws.Range("K134").Value = 19201.5
$ws.Range("M134").Value = -16666.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 14290668
$ws.Range("J68").Value = 16668279
$ws.Range("L68").Value = 50004837
$ws.Range("N68").Value = -50006459
$ws.Range("H71").Value = 14290668
$ws.Range("J71").Value = 16668279
$ws.Range("L71").Value = 150014511
$ws.Range("N71").Value = -150022623
$ws.Range("H80").Value = 3200
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3200
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9600
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -11472
$ws.Range("H83").Value = 3200
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3200
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 28800
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -38160
$ws.Range("H107").Value = 606.36365
$ws.Range("J107").Value = 583.5714
$ws.Range("L107").Value = 1750.7142
$ws.Range("N107").Value = -5590.7142

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 663.5714
$ws.Range("I2").Value = 858.55554
$ws.Range("J2").Value = 312.6
$ws.Range("K2").Value = 858.55554
$ws.Range("L2").Value = 312.6
$ws.Range("M2").Value = -745.55554
$ws.Range("N2").Value = -538.6
$ws.Range("H70").Value = 95723.45
$ws.Range("J70").Value = 4949.4
$ws.Range("L70").Value = 4949.4
$ws.Range("N70").Value = -5489.4
$ws.Range("H73").Value = 95723.45
$ws.Range("J73").Value = 4949.4
$ws.Range("L73").Value = 4949.4
$ws.Range("N73").Value = -6821.4
$ws.Range("H132").Value = 4024.8
$ws.Range("I132").Value = 4369.6
$ws.Range("J132").Value = 3680
$ws.Range("K132").Value = 13108.8
$ws.Range("L132").Value = 11040
$ws.Range("M132").Value = -10578.8
$ws.Range("N132").Value = -16100

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1638.6666
$ws.Range("I61").Value = 1513.8
$ws.Range("K61").Value = 1513.8
$ws.Range("M61").Value = -1311.8
$ws.Range("H68").Value = 4666.6665
$ws.Range("I68").Value = 5500
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 5500
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -4751
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 4666.6665
$ws.Range("I71").Value = 5500
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 27500
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -23756
$ws.Range("N71").Value = -22488
$ws.Range("H113").Value = 1638.6666
$ws.Range("I113").Value = 1513.8
$ws.Range("K113").Value = 1513.8
$ws.Range("M113").Value = 656.2
$ws.Range("H136").Value = 5055.2173
$ws.Range("I136").Value = 5086.2666
$ws.Range("K136").Value = 15258.7998
$ws.Range("M136").Value = -12708.7998

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6449.6
$ws.Range("I81").Value = 6812
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 13624
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -12563
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 6449.6
$ws.Range("I84").Value = 6812
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 68120
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -62816
$ws.Range("N84").Value = -60608
$ws.Range("H122").Value = 13894031
$ws.Range("I122").Value = 5650.5713
$ws.Range("J122").Value = 62503360
$ws.Range("K122").Value = 16951.7139
$ws.Range("L122").Value = 187510080
$ws.Range("M122").Value = -14501.7139
$ws.Range("N122").Value = -187514980
$ws.Range("H132").Value = 11114473
$ws.Range("I132").Value = 15154804
$ws.Range("K132").Value = 45464412
$ws.Range("M132").Value = -45461882
$ws.Range("H136").Value = 27029606
$ws.Range("I136").Value = 32259278
$ws.Range("K136").Value = 96777834
$ws.Range("M136").Value = -96775284
